# AutoCommit-style edit: user scrolled the frozen grid down a bit and
# selected a new cell, then filled in the two trailing grade columns
# (G/H) plus the bonus columns (I/J) for row 30 (Французов Константин).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the active selection in the frozen (bottom-right) pane ---
# (the sheet is already frozen at xSplit=2 / ySplit=3 from C4; just move
# the active cell like the user did when they scrolled down to row 19)
$ws.Range("I19").Select()

# --- Fill in row 30's remaining grade cells ---
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 5

# I30/J30 are brand-new cells in this row; copy the existing formatting
# used by the same columns in the surrounding rows (green fill + thick
# border, style index 6) before writing their values.
$ws.Range("I29").Copy()
$ws.Range("I30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I30").Value = 5

$ws.Range("J29").Copy()
$ws.Range("J30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J30").Value = 5

$excel.CutCopyMode = $false
